$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-1650996074688727"
$wb.Worksheets.Item(2).Name = "NB_TO-16509960768967464"
$wb.Worksheets.Item(3).Name = "RS_TO-16509960768967464"
$wb.Worksheets.Item(4).Name = "TOL_TO-16509960769527245"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16509960770167236"

# Sheet1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-1650996074648688.csv"
$ws1.Range("B3").Value = "GNG_stims-16509960746727293.csv"
$ws1.Range("B4").Value = "go_stims-16509960746727293.csv"
$ws1.Range("B5").Value = "GNG_stims-1650996074688727.csv"

# Sheet2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16509960768086894.csv"
$ws2.Range("B3").Value = "ZB-match_3-16509960757287326.csv"
$ws2.Range("B4").Value = "OB-16509960763847232.csv"
$ws2.Range("B5").Value = "ZB-match_3-16509960748247285.csv"
$ws2.Range("B6").Value = "OB-16509960760086868.csv"
$ws2.Range("B7").Value = "OB-1650996076512692.csv"
$ws2.Range("B8").Value = "TB-16509960768726878.csv"
$ws2.Range("B9").Value = "ZB-match_6-1650996075904722.csv"
$ws2.Range("B10").Value = "TB-16509960767526922.csv"

# Sheet4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16509960769127219.csv"
$ws4.Range("B3").Value = "ZM_stims-16509960768967464.csv"
$ws4.Range("B4").Value = "MM_stims-16509960769367228.csv"
$ws4.Range("B5").Value = "ZM_stims-16509960769127219.csv"
$ws4.Range("B6").Value = "MM_stims-16509960769527245.csv"
$ws4.Range("B7").Value = "ZM_stims-16509960769367228.csv"

# Sheet5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16509960770007365.csv"
$ws5.Range("B3").Value = "SAT_stims-16509960769686978.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509960769847264.csv"
$ws5.Range("B5").Value = "SAT_stims-16509960769527245.csv"
